$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 6700
$ws.Range("J16").Value = 13000
$ws.Range("L16").Value = 13000
$ws.Range("N16").Value = -13460

$ws.Range("H33").Value = 85.75
$ws.Range("I33").Value = 55.666668
$ws.Range("K33").Value = 55.666668
$ws.Range("M33").Value = 173.333332

$ws.Range("H40").Value = 7941
$ws.Range("I40").Value = 7875.375
$ws.Range("J40").Value = 7984.75
$ws.Range("K40").Value = 7875.375
$ws.Range("L40").Value = 7984.75
$ws.Range("M40").Value = -7700.375
$ws.Range("N40").Value = -8334.75

$ws.Range("H107").Value = 33652.633
$ws.Range("I107").Value = 36043.895
$ws.Range("K107").Value = 36043.895
$ws.Range("M107").Value = -34123.895

$ws.Range("H126").Value = 77255.664
$ws.Range("J126").Value = 77255.664
$ws.Range("L126").Value = 77255.664
$ws.Range("N126").Value = -87135.664

$ws.Range("H132").Value = 3330
$ws.Range("I132").Value = 3189.0527
$ws.Range("K132").Value = 9567.158100000001
$ws.Range("M132").Value = -7037.158100000001

$ws.Range("H138").Value = 6748.5273
$ws.Range("I138").Value = 1567.1111
$ws.Range("J138").Value = 9269.216
$ws.Range("K138").Value = 4701.3333
$ws.Range("L138").Value = 27807.648
$ws.Range("M138").Value = 438.6666999999998
$ws.Range("N138").Value = -38087.648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 18400
$ws.Range("J49").Value = 18400
$ws.Range("L49").Value = 18400
$ws.Range("N49").Value = -18920

$ws.Range("H101").Value = 78333.336
$ws.Range("J101").Value = 78333.336
$ws.Range("L101").Value = 78333.336
$ws.Range("N101").Value = -84823.336

$ws.Range("H102").Value = 1521.3334
$ws.Range("I102").Value = 1521.3334
$ws.Range("K102").Value = 1521.3334
$ws.Range("M102").Value = 100.6666

$ws.Range("H122").Value = 4143.4287
$ws.Range("I122").Value = 2336.7693
$ws.Range("K122").Value = 7010.3079
$ws.Range("M122").Value = -4560.3079

$ws.Range("H132").Value = 3870.551
$ws.Range("I132").Value = 2884.8462
$ws.Range("K132").Value = 8654.5386
$ws.Range("M132").Value = -6124.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3226.4285
$ws.Range("I99").Value = 2264.8
$ws.Range("J99").Value = 5630.5
$ws.Range("K99").Value = 2264.8
$ws.Range("L99").Value = 5630.5
$ws.Range("M99").Value = -766.8000000000002
$ws.Range("N99").Value = -8626.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31581.973
$ws.Range("I31").Value = 1764.2174
$ws.Range("J31").Value = 80568.28999999999
$ws.Range("K31").Value = 1764.2174
$ws.Range("L31").Value = 80568.28999999999
$ws.Range("M31").Value = -1469.2174
$ws.Range("N31").Value = -81158.28999999999

$ws.Range("H34").Value = 31581.973
$ws.Range("I34").Value = 1764.2174
$ws.Range("J34").Value = 80568.28999999999
$ws.Range("K34").Value = 1764.2174
$ws.Range("L34").Value = 80568.28999999999
$ws.Range("M34").Value = -1562.2174
$ws.Range("N34").Value = -80972.28999999999

$ws.Range("H58").Value = 6265.263
$ws.Range("I58").Value = 5043.8335
$ws.Range("J58").Value = 8359.143
$ws.Range("K58").Value = 5043.8335
$ws.Range("L58").Value = 8359.143
$ws.Range("M58").Value = -4840.8335
$ws.Range("N58").Value = -8765.143

$ws.Range("H99").Value = 5641.1113
$ws.Range("I99").Value = 4868.1333
$ws.Range("J99").Value = 6607.3335
$ws.Range("K99").Value = 4868.1333
$ws.Range("L99").Value = 6607.3335
$ws.Range("M99").Value = -3370.1333
$ws.Range("N99").Value = -9603.333500000001

$ws.Range("H126").Value = 5641.1113
$ws.Range("I126").Value = 4868.1333
$ws.Range("J126").Value = 6607.3335
$ws.Range("K126").Value = 14604.3999
$ws.Range("L126").Value = 19822.0005
$ws.Range("M126").Value = -12134.3999
$ws.Range("N126").Value = -24762.0005

$ws.Range("H130").Value = 56641.75
$ws.Range("J130").Value = 56641.75
$ws.Range("L130").Value = 56641.75
$ws.Range("N130").Value = -66681.75

$ws.Range("H132").Value = 6012.3887
$ws.Range("I132").Value = 5881.1665
$ws.Range("J132").Value = 6274.8335
$ws.Range("K132").Value = 17643.4995
$ws.Range("L132").Value = 18824.5005
$ws.Range("M132").Value = -15113.4995
$ws.Range("N132").Value = -23884.5005

$ws.Range("H136").Value = 6265.263
$ws.Range("I136").Value = 5043.8335
$ws.Range("J136").Value = 8359.143
$ws.Range("K136").Value = 15131.5005
$ws.Range("L136").Value = 25077.429
$ws.Range("M136").Value = -12581.5005
$ws.Range("N136").Value = -30177.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 408.66666
$ws.Range("J2").Value = 750
$ws.Range("L2").Value = 4500
$ws.Range("N2").Value = -4726

$ws.Range("H5").Value = 635674.75
$ws.Range("I5").Value = 57461.855
$ws.Range("K5").Value = 172385.565
$ws.Range("M5").Value = -172273.565

$ws.Range("H92").Value = 556348.9399999999
$ws.Range("I92").Value = 769752.4
$ws.Range("J92").Value = 1500
$ws.Range("K92").Value = 2309257.2
$ws.Range("L92").Value = 4500
$ws.Range("M92").Value = -2308009.2
$ws.Range("N92").Value = -6996

$ws.Range("H107").Value = 48634.363
$ws.Range("I107").Value = 979.9167
$ws.Range("J107").Value = 105819.7
$ws.Range("K107").Value = 2939.7501
$ws.Range("L107").Value = 317459.1
$ws.Range("M107").Value = -1019.7501
$ws.Range("N107").Value = -321299.1

$ws.Range("H131").Value = 11566415
$ws.Range("I131").Value = 27861732
$ws.Range("J131").Value = 63839.06
$ws.Range("K131").Value = 83585196
$ws.Range("L131").Value = 191517.18
$ws.Range("M131").Value = -83580156
$ws.Range("N131").Value = -201597.18

$ws.Range("H135").Value = 635674.75
$ws.Range("I135").Value = 57461.855
$ws.Range("K135").Value = 517156.695
$ws.Range("M135").Value = -514621.695

$ws.Range("H137").Value = 2052.6667
$ws.Range("I137").Value = 1336.75
$ws.Range("K137").Value = 4010.25
$ws.Range("M137").Value = 1089.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 19310
$ws.Range("J47").Value = 19310
$ws.Range("L47").Value = 19310
$ws.Range("N47").Value = -20446

$ws.Range("H55").Value = 18378.75
$ws.Range("I55").Value = 18676.666
$ws.Range("K55").Value = 18676.666
$ws.Range("M55").Value = -18349.666

$ws.Range("H80").Value = 1433057.2
$ws.Range("I80").Value = 1114921.9
$ws.Range("J80").Value = 2005701
$ws.Range("K80").Value = 1114921.9
$ws.Range("L80").Value = 2005701
$ws.Range("M80").Value = -1113923.9
$ws.Range("N80").Value = -2007697

$ws.Range("H83").Value = 1433057.2
$ws.Range("I83").Value = 1114921.9
$ws.Range("J83").Value = 2005701
$ws.Range("K83").Value = 5574609.5
$ws.Range("L83").Value = 10028505
$ws.Range("M83").Value = -5569617.5
$ws.Range("N83").Value = -10038489

$ws.Range("H122").Value = 3472.75
$ws.Range("I122").Value = 2046.1818
$ws.Range("K122").Value = 6138.5454
$ws.Range("M122").Value = -3688.5454

$ws.Range("H126").Value = 3897.75
$ws.Range("I126").Value = 2587
$ws.Range("K126").Value = 7761
$ws.Range("M126").Value = -5291

$ws.Range("H132").Value = 326467.7
$ws.Range("I132").Value = 419841.72
$ws.Range("J132").Value = 102370
$ws.Range("K132").Value = 1259525.16
$ws.Range("L132").Value = 307110
$ws.Range("M132").Value = -1256995.16
$ws.Range("N132").Value = -312170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 17370
$ws.Range("J41").Value = 17370
$ws.Range("L41").Value = 17370
$ws.Range("N41").Value = -18246

$ws.Range("H46").Value = 5456.7617
$ws.Range("I46").Value = 5649.4287
$ws.Range("K46").Value = 5649.4287
$ws.Range("M46").Value = -5461.4287

$ws.Range("H47").Value = 18650
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H52").Value = 18650
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H132").Value = 3866
$ws.Range("I132").Value = 2407.25
$ws.Range("K132").Value = 7221.75
$ws.Range("M132").Value = -4691.75

$ws.Range("H139").Value = 49800
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49800
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49800
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -60080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 16866.334
$ws.Range("I3").Value = 49999
$ws.Range("K3").Value = 49999
$ws.Range("M3").Value = -49885

$ws.Range("H96").Value = 143877.72
$ws.Range("I96").Value = 333988.66
$ws.Range("K96").Value = 333988.66
$ws.Range("M96").Value = -332615.66

$ws.Range("H122").Value = 43483396
$ws.Range("I122").Value = 55560230
$ws.Range("K122").Value = 166680690
$ws.Range("M122").Value = -166678240

$ws.Range("H139").Value = 88888
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 88888
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 88888
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -99168
